# Weekly update: insert a new price record (row 266) for Melón / Tuna / Primera
# at Vega Monumental Concepción, pushing the existing rows 266-303 down to
# 267-304.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 266 - shifts rows 266..303 down to 267..304
# and extends the used range to A1:R304.
$ws.Rows.Item(266).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A266").Value = 11
$ws.Range("B266").Value = "Vega Monumental Concepción"
$ws.Range("C266").Value = "Bíobío"
$ws.Range("D266").Value = 44932
$ws.Range("E266").Value = 8
$ws.Range("F266").Value = 100112027
$ws.Range("G266").Value = "Melón"
$ws.Range("H266").Value = "Tuna"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 2200
$ws.Range("K266").Value = 800
$ws.Range("L266").Value = 900
$ws.Range("M266").Value = 845
$ws.Range("N266").Value = "$/unidad"
$ws.Range("O266").Value = "Región de O'Higgins"
$ws.Range("P266").Value = 845
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = "Hortaliza"
